# "Generate Report for Handoff" - refresh the localization-status report with
# a new handoff run: new source GUID file, new handoff commit hash / xliff
# names, new handoff timestamps, and reset (cleared) handback state since the
# freshly generated targets have not been handed back yet.

$wb = $excel.ActiveWorkbook

$oldGuid = "b9b8b0b6-6e53-45b8-bf1b-cdc454d0f50a"
$newGuid = "092590bd-e68f-4bf5-9d0f-29f8c6a4f22f"

$oldCommit = "010bfedd6e8ea978eb1e60f9dfd1536063fdd740"
$newCommit = "12f9ad841acd8306c71d28750a2a1ef99e0b6bf5"

$newHoDate      = "2016-08-25 21:00:56"
$newZhHoDate    = "2016-08-25 21:00:34"
$newDeHoDate    = "2016-08-25 21:00:56"
$epoch          = "0001-01-01 00:00:00"

$newFileName    = "$newGuid.md"
$newPathName    = "e2e\$newGuid.md"
$newZhXlf       = "$newGuid.$newCommit.zh-cn.xlf"
$newDeXlf       = "$newGuid.$newCommit.de-de.xlf"

$ghUrlBase      = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/92d2561af4499624546d3a472fa293a23b88d361/e2e/$oldGuid.md"

# OOXML stores col width as {ColumnWidth}+5/MaxDigitWidth(~6px) padding; back
# out that constant so the written <col width="..."/> lands on the target.
$padding  = 0.8333333333333333
$colAWidth = 39.7459280831473 - $padding
$colIWidth = 18.6506053379604 - $padding
$colJWidth = 21.7054770333426 - $padding

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Cells.Item(2,1).Value = $newFileName          # A2 File Name
$ov.Cells.Item(2,2).Value = $newPathName          # B2 Path And Name
$ov.Cells.Item(2,7).Value = $newHoDate            # G2 Latest HO Xliff Generate Date

# refresh hyperlink display text on B2 (target URL is unchanged)
$ov.Range("B2").Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("B2"), $ghUrlBase, "", "", $newPathName)

$ov.Columns.Item(1).ColumnWidth = $colAWidth

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Cells.Item(2,1).Value = $newFileName          # A2 Source File Name
$zh.Cells.Item(2,7).Value = $newZhXlf             # G2 Latest Handoff File
$zh.Cells.Item(2,8).Value = $newZhHoDate          # H2 Latest Handoff Datetime
$zh.Cells.Item(2,9).Value = ""                    # I2 Latest Target File (cleared)
$zh.Cells.Item(2,9).Style = "Normal"
$zh.Cells.Item(2,10).Value = ""                   # J2 Latest Handback File (cleared)
$zh.Cells.Item(2,11).Value = $epoch               # K2 Latest Handback DateTime (reset)

# refresh hyperlink display text on A2, drop the I2 (target file) hyperlink
$zh.Range("A2").Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), $ghUrlBase, "", "", $newFileName)

$zh.Columns.Item(1).ColumnWidth = $colAWidth
$zh.Columns.Item(9).ColumnWidth = $colIWidth
$zh.Columns.Item(10).ColumnWidth = $colJWidth

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Cells.Item(2,1).Value = $newFileName          # A2 Source File Name
$de.Cells.Item(2,7).Value = $newDeXlf             # G2 Latest Handoff File
$de.Cells.Item(2,8).Value = $newDeHoDate          # H2 Latest Handoff Datetime
$de.Cells.Item(2,9).Value = ""                    # I2 Latest Target File (cleared)
$de.Cells.Item(2,9).Style = "Normal"
$de.Cells.Item(2,10).Value = ""                   # J2 Latest Handback File (cleared)
$de.Cells.Item(2,11).Value = $epoch               # K2 Latest Handback DateTime (reset)

# refresh hyperlink display text on A2, drop the I2 (target file) hyperlink
$de.Range("A2").Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), $ghUrlBase, "", "", $newFileName)

$de.Columns.Item(1).ColumnWidth = $colAWidth
$de.Columns.Item(9).ColumnWidth = $colIWidth
$de.Columns.Item(10).ColumnWidth = $colJWidth
